$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Insert a new row at the top, shifting the existing table down by one row.
$ws.Rows.Item(1).Insert()

# New small summary table in columns I:K, rows 1-3.
$ws.Range("A1").Value = "Trial 1"
$ws.Range("I1").Value = "Trial 2"

$ws.Range("I2").Value = "Left"
$ws.Range("J2").Value = "Right"
$ws.Range("K2").Value = "Straight"

$ws.Range("I3").Value = 734
$ws.Range("J3").Value = 323
$ws.Range("K3").Value = 243

# Match the recorded selection in the saved file.
$ws.Range("I4").Select()
